$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 545.6739  # H17: 550.0217 -> 545.6739
$ws.Cells.Item(17, 10).Value = 545.6739  # J17: 550.0217 -> 545.6739
$ws.Cells.Item(17, 12).Value = 1637.0217  # L17: 1650.0651 -> 1637.0217
$ws.Cells.Item(17, 14).Value = -1973.0217  # N17: -1986.0651 -> -1973.0217
$ws.Cells.Item(132, 8).Value = 10106590  # H132: 12827467 -> 10106590
$ws.Cells.Item(132, 9).Value = 13895579  # I132: 17552238 -> 13895579
$ws.Cells.Item(132, 10).Value = 2616.6667  # J132: 3085.8572 -> 2616.6667
$ws.Cells.Item(132, 11).Value = 41686737  # K132: 52656714 -> 41686737
$ws.Cells.Item(132, 12).Value = 7850.000100000001  # L132: 9257.571599999999 -> 7850.000100000001
$ws.Cells.Item(132, 13).Value = -41684207  # M132: -52654184 -> -41684207
$ws.Cells.Item(132, 14).Value = -12910.0001  # N132: -14317.5716 -> -12910.0001
$ws.Cells.Item(137, 8).Value = 1130.2  # H137: 1168.6052 -> 1130.2
$ws.Cells.Item(137, 9).Value = 769.8889  # I137: 799.4400000000001 -> 769.8889
$ws.Cells.Item(137, 11).Value = 2309.6667  # K137: 2398.32 -> 2309.6667
$ws.Cells.Item(137, 13).Value = 240.3332999999998  # M137: 151.6799999999998 -> 240.3332999999998
$ws.Cells.Item(138, 8).Value = 1516.83  # H138: 454558.97 -> 1516.83
$ws.Cells.Item(138, 9).Value = 658.63635  # I138: 1097.2903 -> 658.63635
$ws.Cells.Item(138, 10).Value = 1939.5223  # J138: 670825.3 -> 1939.5223
$ws.Cells.Item(138, 11).Value = 1975.90905  # K138: 3291.8709 -> 1975.90905
$ws.Cells.Item(138, 12).Value = 5818.5669  # L138: 2012475.9 -> 5818.5669
$ws.Cells.Item(138, 13).Value = 3164.09095  # M138: 1848.1291 -> 3164.09095
$ws.Cells.Item(138, 14).Value = -16098.5669  # N138: -2022755.9 -> -16098.5669
$ws.Cells.Item(141, 8).Value = 909.5  # H141: 1072.7778 -> 909.5
$ws.Cells.Item(141, 9).Value = 710.25  # I141: 777.8570999999999 -> 710.25
$ws.Cells.Item(141, 11).Value = 2130.75  # K141: 2333.5713 -> 2130.75
$ws.Cells.Item(141, 13).Value = 3049.25  # M141: 2846.4287 -> 3049.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3584.6086  # H32: 3431.959 -> 3584.6086
$ws.Cells.Item(32, 9).Value = 3166.5322  # I32: 3023.0303 -> 3166.5322
$ws.Cells.Item(32, 11).Value = 3166.5322  # K32: 3023.0303 -> 3166.5322
$ws.Cells.Item(32, 13).Value = -2879.5322  # M32: -2736.0303 -> -2879.5322
$ws.Cells.Item(61, 8).Value = 1377.0625  # H61: 1445 -> 1377.0625
$ws.Cells.Item(61, 9).Value = 1230.88  # I61: 1294.5834 -> 1230.88
$ws.Cells.Item(61, 10).Value = 1899.1428  # J61: 2046.6666 -> 1899.1428
$ws.Cells.Item(61, 11).Value = 1230.88  # K61: 1294.5834 -> 1230.88
$ws.Cells.Item(61, 12).Value = 1899.1428  # L61: 2046.6666 -> 1899.1428
$ws.Cells.Item(61, 13).Value = -1018.88  # M61: -1082.5834 -> -1018.88
$ws.Cells.Item(61, 14).Value = -2323.1428  # N61: -2470.6666 -> -2323.1428
$ws.Cells.Item(74, 8).Value = 993.85187  # H74: 1165.0952 -> 993.85187
$ws.Cells.Item(74, 9).Value = 740.6087  # I74: 833.35297 -> 740.6087
$ws.Cells.Item(74, 10).Value = 2450  # J74: 2575 -> 2450
$ws.Cells.Item(74, 11).Value = 740.6087  # K74: 833.35297 -> 740.6087
$ws.Cells.Item(74, 12).Value = 2450  # L74: 2575 -> 2450
$ws.Cells.Item(74, 13).Value = 133.3913  # M74: 40.64702999999997 -> 133.3913
$ws.Cells.Item(74, 14).Value = -4198  # N74: -4323 -> -4198
$ws.Cells.Item(77, 8).Value = 993.85187  # H77: 1165.0952 -> 993.85187
$ws.Cells.Item(77, 9).Value = 740.6087  # I77: 833.35297 -> 740.6087
$ws.Cells.Item(77, 10).Value = 2450  # J77: 2575 -> 2450
$ws.Cells.Item(77, 11).Value = 3703.0435  # K77: 4166.76485 -> 3703.0435
$ws.Cells.Item(77, 12).Value = 12250  # L77: 12875 -> 12250
$ws.Cells.Item(77, 13).Value = 664.9565000000002  # M77: 201.2351499999995 -> 664.9565000000002
$ws.Cells.Item(77, 14).Value = -20986  # N77: -21611 -> -20986
$ws.Cells.Item(88, 8).Value = 3023.4614  # H88: 2986.0715 -> 3023.4614
$ws.Cells.Item(88, 10).Value = 3127.2727  # J88: 3075 -> 3127.2727
$ws.Cells.Item(88, 12).Value = 3127.2727  # L88: 3075 -> 3127.2727
$ws.Cells.Item(88, 14).Value = -3939.2727  # N88: -3887 -> -3939.2727
$ws.Cells.Item(91, 8).Value = 3023.4614  # H91: 2986.0715 -> 3023.4614
$ws.Cells.Item(91, 10).Value = 3127.2727  # J91: 3075 -> 3127.2727
$ws.Cells.Item(91, 12).Value = 3127.2727  # L91: 3075 -> 3127.2727
$ws.Cells.Item(91, 14).Value = -5935.2727  # N91: -5883 -> -5935.2727
$ws.Cells.Item(132, 8).Value = 2572.1333  # H132: 2394.4375 -> 2572.1333
$ws.Cells.Item(132, 9).Value = 2234.9092  # I132: 2126.0833 -> 2234.9092
$ws.Cells.Item(132, 10).Value = 3499.5  # J132: 3199.5 -> 3499.5
$ws.Cells.Item(132, 11).Value = 6704.7276  # K132: 6378.249899999999 -> 6704.7276
$ws.Cells.Item(132, 12).Value = 10498.5  # L132: 9598.5 -> 10498.5
$ws.Cells.Item(132, 13).Value = -4174.7276  # M132: -3848.249899999999 -> -4174.7276
$ws.Cells.Item(132, 14).Value = -15558.5  # N132: -14658.5 -> -15558.5
$ws.Cells.Item(134, 8).Value = 34999.57  # H134: 35000 -> 34999.57
$ws.Cells.Item(134, 10).Value = 34999.57  # J134: 35000 -> 34999.57
$ws.Cells.Item(134, 12).Value = 34999.57  # L134: 35000 -> 34999.57
$ws.Cells.Item(134, 14).Value = -45139.57  # N134: -45140 -> -45139.57
$ws.Cells.Item(136, 8).Value = 1377.0625  # H136: 1445 -> 1377.0625
$ws.Cells.Item(136, 9).Value = 1230.88  # I136: 1294.5834 -> 1230.88
$ws.Cells.Item(136, 10).Value = 1899.1428  # J136: 2046.6666 -> 1899.1428
$ws.Cells.Item(136, 11).Value = 3692.64  # K136: 3883.7502 -> 3692.64
$ws.Cells.Item(136, 12).Value = 5697.428400000001  # L136: 6139.9998 -> 5697.428400000001
$ws.Cells.Item(136, 13).Value = -1142.64  # M136: -1333.7502 -> -1142.64
$ws.Cells.Item(136, 14).Value = -10797.4284  # N136: -11239.9998 -> -10797.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4968.2104  # H86: 5038.722 -> 4968.2104
$ws.Cells.Item(86, 9).Value = 5087.3125  # I86: 5179.8667 -> 5087.3125
$ws.Cells.Item(86, 11).Value = 5087.3125  # K86: 5179.8667 -> 5087.3125
$ws.Cells.Item(86, 13).Value = -3964.3125  # M86: -4056.8667 -> -3964.3125
$ws.Cells.Item(89, 8).Value = 4968.2104  # H89: 5038.722 -> 4968.2104
$ws.Cells.Item(89, 9).Value = 5087.3125  # I89: 5179.8667 -> 5087.3125
$ws.Cells.Item(89, 11).Value = 25436.5625  # K89: 25899.3335 -> 25436.5625
$ws.Cells.Item(89, 13).Value = -19820.5625  # M89: -20283.3335 -> -19820.5625
$ws.Cells.Item(105, 8).Value = 66669540  # H105: 62502944 -> 66669540
$ws.Cells.Item(105, 9).Value = 71431496  # I105: 66669660 -> 71431496
$ws.Cells.Item(105, 11).Value = 71431496  # K105: 66669660 -> 71431496
$ws.Cells.Item(105, 13).Value = -71429749  # M105: -66667913 -> -71429749

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1157.0476  # H31: 1216.069 -> 1157.0476
$ws.Cells.Item(31, 9).Value = 1134.4642  # I31: 1199.3726 -> 1134.4642
$ws.Cells.Item(31, 11).Value = 1134.4642  # K31: 1199.3726 -> 1134.4642
$ws.Cells.Item(31, 13).Value = -839.4641999999999  # M31: -904.3725999999999 -> -839.4641999999999
$ws.Cells.Item(34, 8).Value = 1157.0476  # H34: 1216.069 -> 1157.0476
$ws.Cells.Item(34, 9).Value = 1134.4642  # I34: 1199.3726 -> 1134.4642
$ws.Cells.Item(34, 11).Value = 1134.4642  # K34: 1199.3726 -> 1134.4642
$ws.Cells.Item(34, 13).Value = -932.4641999999999  # M34: -997.3725999999999 -> -932.4641999999999
$ws.Cells.Item(58, 8).Value = 855.2759  # H58: 876.86206 -> 855.2759
$ws.Cells.Item(58, 9).Value = 794.9545000000001  # I58: 823.4091 -> 794.9545000000001
$ws.Cells.Item(58, 11).Value = 794.9545000000001  # K58: 823.4091 -> 794.9545000000001
$ws.Cells.Item(58, 13).Value = -591.9545000000001  # M58: -620.4091 -> -591.9545000000001
$ws.Cells.Item(134, 8).Value = 1098.2632  # H134: 2636.1333 -> 1098.2632
$ws.Cells.Item(134, 9).Value = 883.8182  # I134: 3260.889 -> 883.8182
$ws.Cells.Item(134, 10).Value = 1393.125  # J134: 1699 -> 1393.125
$ws.Cells.Item(134, 11).Value = 2651.4546  # K134: 9782.667000000001 -> 2651.4546
$ws.Cells.Item(134, 12).Value = 4179.375  # L134: 5097 -> 4179.375
$ws.Cells.Item(134, 13).Value = -116.4546  # M134: -7247.667000000001 -> -116.4546
$ws.Cells.Item(134, 14).Value = -9249.375  # N134: -10167 -> -9249.375
$ws.Cells.Item(136, 8).Value = 855.2759  # H136: 876.86206 -> 855.2759
$ws.Cells.Item(136, 9).Value = 794.9545000000001  # I136: 823.4091 -> 794.9545000000001
$ws.Cells.Item(136, 11).Value = 2384.8635  # K136: 2470.2273 -> 2384.8635
$ws.Cells.Item(136, 13).Value = 165.1364999999996  # M136: 79.77269999999999 -> 165.1364999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 2454.5454  # H55: 2666.6667 -> 2454.5454
$ws.Cells.Item(55, 10).Value = 2777.7778  # J55: 2875 -> 2777.7778
$ws.Cells.Item(55, 12).Value = 8333.3334  # L55: 8625 -> 8333.3334
$ws.Cells.Item(55, 14).Value = -8687.3334  # N55: -8979 -> -8687.3334
$ws.Cells.Item(107, 8).Value = 8772.166999999999  # H107: 9552.362999999999 -> 8772.166999999999
$ws.Cells.Item(107, 9).Value = 407.5  # I107: 480 -> 407.5
$ws.Cells.Item(107, 11).Value = 1222.5  # K107: 1440 -> 1222.5
$ws.Cells.Item(107, 13).Value = 697.5  # M107: 480 -> 697.5
$ws.Cells.Item(113, 8).Value = 632.4054  # H113: 632.8684 -> 632.4054
$ws.Cells.Item(113, 10).Value = 654.67645  # J113: 654.54285 -> 654.67645
$ws.Cells.Item(113, 12).Value = 1964.02935  # L113: 1963.62855 -> 1964.02935
$ws.Cells.Item(113, 14).Value = -6304.029350000001  # N113: -6303.62855 -> -6304.029350000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 10325.385  # H123: 10326 -> 10325.385
$ws.Cells.Item(123, 10).Value = 10325.385  # J123: 10326 -> 10325.385
$ws.Cells.Item(123, 12).Value = 10325.385  # L123: 10326 -> 10325.385
$ws.Cells.Item(123, 14).Value = -15225.385  # N123: -15226 -> -15225.385
$ws.Cells.Item(132, 8).Value = 2089.4595  # H132: 2300.5 -> 2089.4595
$ws.Cells.Item(132, 9).Value = 1495.8636  # I132: 1874.6522 -> 1495.8636
$ws.Cells.Item(132, 10).Value = 2960.0667  # J132: 2953.4666 -> 2960.0667
$ws.Cells.Item(132, 11).Value = 4487.5908  # K132: 5623.9566 -> 4487.5908
$ws.Cells.Item(132, 12).Value = 8880.2001  # L132: 8860.399800000001 -> 8880.2001
$ws.Cells.Item(132, 13).Value = -1957.5908  # M132: -3093.9566 -> -1957.5908
$ws.Cells.Item(132, 14).Value = -13940.2001  # N132: -13920.3998 -> -13940.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 0  # H100: 2000 -> 0
$ws.Cells.Item(100, 10).Value = 0  # J100: 2000 -> 0
$ws.Cells.Item(100, 12).Value = 0  # L100: 2000 -> 0
$ws.Cells.Item(100, 14).ClearContents()  # N100 removed (was -3082)
$ws.Cells.Item(132, 8).Value = 21065.588  # H132: 21484.06 -> 21065.588
$ws.Cells.Item(132, 9).Value = 979.19354  # I132: 1007.1 -> 979.19354
$ws.Cells.Item(132, 11).Value = 2937.58062  # K132: 3021.3 -> 2937.58062
$ws.Cells.Item(132, 13).Value = -407.5806199999997  # M132: -491.3000000000002 -> -407.5806199999997
$ws.Cells.Item(136, 8).Value = 974.2593000000001  # H136: 1073.7916 -> 974.2593000000001
$ws.Cells.Item(136, 9).Value = 841.6667  # I136: 936.4761999999999 -> 841.6667
$ws.Cells.Item(136, 11).Value = 2525.0001  # K136: 2809.4286 -> 2525.0001
$ws.Cells.Item(136, 13).Value = 24.9998999999998  # M136: -259.4285999999997 -> 24.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value = 0  # H104: 12345 -> 0
$ws.Cells.Item(104, 10).Value = 0  # J104: 12345 -> 0
$ws.Cells.Item(104, 12).Value = 0  # L104: 12345 -> 0
$ws.Cells.Item(104, 14).ClearContents()  # N104 removed (was -19333)
$ws.Cells.Item(122, 8).Value = 11819692  # H122: 8667490 -> 11819692
$ws.Cells.Item(122, 9).Value = 13685707  # I122: 11305204 -> 13685707
$ws.Cells.Item(122, 10).Value = 1600  # J122: 714.1429000000001 -> 1600
$ws.Cells.Item(122, 11).Value = 41057121  # K122: 33915612 -> 41057121
$ws.Cells.Item(122, 12).Value = 4800  # L122: 2142.4287 -> 4800
$ws.Cells.Item(122, 13).Value = -41054671  # M122: -33913162 -> -41054671
$ws.Cells.Item(122, 14).Value = -9700  # N122: -7042.4287 -> -9700
$ws.Cells.Item(132, 8).Value = 2417  # H132: 2438.9285 -> 2417
$ws.Cells.Item(132, 9).Value = 2619.7188  # I132: 2632.875 -> 2619.7188
$ws.Cells.Item(132, 10).Value = 1768.3  # J132: 1818.3 -> 1768.3
$ws.Cells.Item(132, 11).Value = 7859.1564  # K132: 7898.625 -> 7859.1564
$ws.Cells.Item(132, 12).Value = 5304.9  # L132: 5454.9 -> 5304.9
$ws.Cells.Item(132, 13).Value = -5329.1564  # M132: -5368.625 -> -5329.1564
$ws.Cells.Item(132, 14).Value = -10364.9  # N132: -10514.9 -> -10364.9

